$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" sheet, so
#    the final sheet order is 2020-Q4, 2021-Q4, 2022-Q1, 总计.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Add($wsTotal)
$wsQ1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts everything after it over by one
# position, and worksheet handles in this object model are resolved
# by position rather than stable identity - so $wsTotal would now
# silently point at the freshly-added "2022-Q1" sheet unless we
# re-fetch it by name.
$wsTotal = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# 2. Populate the header row of "2022-Q1", copying the formatting
#    (bold / border / center) from the equivalent header on "2021-Q4".
# ------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ1.Cells.Item(1,2).Value2 = "基金代码"
$wsQ1.Cells.Item(1,3).Value2 = "基金名称"
$wsQ1.Cells.Item(1,4).Value2 = "基金规模"
$wsQ1.Cells.Item(1,5).Value2 = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value2 = "仓位占比"
$wsQ1.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value2 = "仓位排名"

# Same index-column styling (bold/border/center) used on A2 of the
# other per-quarter sheets.
$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A4").PasteSpecial(-4122)

# The B (fund code, has leading zeros) and D:G columns hold
# numeric-looking values that are stored as text in the source data,
# so mark the ranges as Text before writing so the values round-trip
# as strings rather than being coerced to numbers.
$wsQ1.Range("B2:B4").NumberFormat = "@"
$wsQ1.Range("D2:G4").NumberFormat = "@"

# ------------------------------------------------------------------
# 3. Data rows for 2022-Q1.
# ------------------------------------------------------------------
$wsQ1.Cells.Item(2,1).Value2 = 0
$wsQ1.Cells.Item(2,2).Value2 = "004532"
$wsQ1.Cells.Item(2,3).Value2 = "民生加银中证港股通高股息精选指数A"
$wsQ1.Cells.Item(2,4).Value2 = "0.26"
$wsQ1.Cells.Item(2,5).Value2 = "94.88"
$wsQ1.Cells.Item(2,6).Value2 = "3.30"
$wsQ1.Cells.Item(2,7).Value2 = "0.0086"
$wsQ1.Cells.Item(2,8).Value2 = 10

$wsQ1.Cells.Item(3,1).Value2 = 1
$wsQ1.Cells.Item(3,2).Value2 = "004533"
$wsQ1.Cells.Item(3,3).Value2 = "民生加银中证港股通高股息精选指数C"
$wsQ1.Cells.Item(3,4).Value2 = "0.10"
$wsQ1.Cells.Item(3,5).Value2 = "94.88"
$wsQ1.Cells.Item(3,6).Value2 = "3.30"
$wsQ1.Cells.Item(3,7).Value2 = "0.0033"
$wsQ1.Cells.Item(3,8).Value2 = 10

$wsQ1.Cells.Item(4,1).Value2 = 2
$wsQ1.Cells.Item(4,2).Value2 = "005770"
$wsQ1.Cells.Item(4,3).Value2 = "信达澳银中证沪港深高股息精选指数"
$wsQ1.Cells.Item(4,4).Value2 = "0.01"
$wsQ1.Cells.Item(4,5).Value2 = "92.47"
$wsQ1.Cells.Item(4,6).Value2 = "2.51"
$wsQ1.Cells.Item(4,7).Value2 = "0.0003"
$wsQ1.Cells.Item(4,8).Value2 = 6

# ------------------------------------------------------------------
# 4. Update the "总计" (totals) sheet: push its two existing rows down
#    by one and add a new leading row summarising 2022-Q1.
# ------------------------------------------------------------------

# old row 3 (2020-Q4) -> row 4
$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Cells.Item(4,1).PasteSpecial(-4122)
$wsTotal.Cells.Item(4,1).Value2 = 2
$wsTotal.Cells.Item(4,2).Value2 = $wsTotal.Cells.Item(3,2).Value2
$wsTotal.Cells.Item(4,3).Value2 = $wsTotal.Cells.Item(3,3).Value2
$wsTotal.Cells.Item(4,4).Value2 = $wsTotal.Cells.Item(3,4).Value2

# old row 2 (2021-Q4) -> row 3
$wsTotal.Cells.Item(2,1).Copy()
$wsTotal.Cells.Item(3,1).PasteSpecial(-4122)
$wsTotal.Cells.Item(3,1).Value2 = 1
$wsTotal.Cells.Item(3,2).Value2 = $wsTotal.Cells.Item(2,2).Value2
$wsTotal.Cells.Item(3,3).Value2 = $wsTotal.Cells.Item(2,3).Value2
$wsTotal.Cells.Item(3,4).Value2 = $wsTotal.Cells.Item(2,4).Value2

# new row 2 (2022-Q1)
$wsTotal.Cells.Item(2,2).Value2 = "2022-Q1"
$wsTotal.Cells.Item(2,3).Value2 = 3
$wsTotal.Cells.Item(2,4).Value2 = 0.01
